$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Periodo Mora" (E) / "Valor Mora" (F) / "Salario Basico" (G) table ---
# The batch of account-statement periods (rows 16-32) is being reversed in order
# (oldest period first instead of newest first) as part of adding the new
# "parte 1" periods, and the Salario Basico (G) is refreshed to its new value.

$periods = @("1708","1709","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903")
$valorMora = @(27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,25740)
$salarioBasico = 689455

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value2 = $periods[$i]
    $ws.Cells.Item($row, 6).Value2 = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value2 = $salarioBasico
}

# --- Reposition the logo image slightly to the left ---
# (Width is re-asserted to its current value so the shape's absolute
# EMU geometry is recomputed precisely alongside the new Left.)
$shp = $ws.Shapes.Item(1)
$shp.Width = 76.81889763779527
$shp.Left = 50.84055118110236
